# [Fonds de solidarite] Add 2020-12-31 data
# Updates nombre_aides (col C) and montant_total (col D) for a handful of
# region/classe_effectif rows, preserving the original text-cell data type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$edits = @(
  @{ Row = 10;  C = "482";  D = "2631981.26" },
  @{ Row = 12;  C = "65";   D = "873168.55" },
  @{ Row = 16;  C = "530";  D = "2687400.85" },
  @{ Row = 17;  C = "243";  D = "2563279.45" },
  @{ Row = 60;  C = "15";   D = "363767.27" },
  @{ Row = 90;  C = "289";  D = "1308853.14" },
  @{ Row = 91;  C = "1103"; D = "5860114.39" },
  @{ Row = 105; C = "495";  D = "2306017.74" },
  @{ Row = 106; C = "231";  D = "1756097.68" },
  @{ Row = 108; C = "18";   D = "357957.62" },
  @{ Row = 112; C = "737";  D = "5120130.94" },
  @{ Row = 114; C = "84";   D = "1250060.00" }
)

foreach ($edit in $edits) {
  # The source data stores every cell as text (inline string), even the
  # numeric-looking ones. Force the NumberFormat to Text before assigning
  # so the new value is kept verbatim as a string instead of being
  # reinterpreted as a number (which would drop things like trailing
  # zeros). Afterwards restore the cell style so formatting stays as it
  # was originally (General / default style).
  $cCell = $ws.Cells.Item($edit.Row, 3)
  $cCell.NumberFormat = "@"
  $cCell.Value = $edit.C
  $cCell.Style = "Normal"

  $dCell = $ws.Cells.Item($edit.Row, 4)
  $dCell.NumberFormat = "@"
  $dCell.Value = $edit.D
  $dCell.Style = "Normal"
}
